$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Fix typo in TC007 summary
$ws.Range("D8").Value = "To verify if the rent car is still available"

# Add new test case row 9 (TC008)
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "TS003"
$ws.Range("C9").Value = "TC008"
$ws.Range("D9").Value = "To verify if the rent car is avaible between two different rent period"
$ws.Range("E9").Value = "RS001"
$ws.Range("F9").Value = "Web is launch"
$ws.Range("G9").Value = "1. Login`n2. Select rent date`n3. Select return date`n4. Press ""Search"" button`n"
$ws.Range("H9").Value = "period1:`nrent date: 20191210`nreturn date: 20191215`nperiod2:`nrent date: 20191220`nreturn date:20191225`nnew rent period:`nrent date:20191216`nreturn date:20191219"
$ws.Range("I9").Value = "Show 4 cars in the window"
$ws.Range("J9").Value = "Only show 3 cars"
$ws.Range("K9").Value = "Failed"
$ws.Range("L9").Value = "Require to fix the defect`nUpdate:`nFixed the defect by introducing a new table rentRecord to keep the rental history"
$ws.Range("M9").Value = "Jimmy Luo"
$ws.Range("N9").Value = [DateTime]::new(2019,12,22)
$ws.Range("O9").Value = "Jimmy Luo"
$ws.Range("P9").Value = [DateTime]::new(2019,12,22)
$ws.Range("Q9").Value = "Browser: Chrome"

Write-Output "done"
